$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column I (Antal) to stay text-typed before writing numeric-looking values,
# matching the source data (values like "1","2","10" stored as text, not numbers).
$ws.Range("I2:I9").NumberFormat = "@"

# Row 2  (was row 5)
$ws.Range("A2").Value2 = 110282835
$ws.Range("B2").Value2 = 89405
$ws.Range("E2").Value2 = 1202
$ws.Range("Q2").Value2 = 600805.3583702671
$ws.Range("R2").Value2 = 6613969.910894822
$ws.Range("D2").Value2 = 'NT'
$ws.Range("F2").Value2 = 'Ullticka'
$ws.Range("G2").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H2").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I2").Value2 = '1'
$ws.Range("J2").Value2 = 'mycel'
$ws.Range("K2").Value2 = ''

# Row 3  (was row 8)
$ws.Range("A3").Value2 = 110282764
$ws.Range("B3").Value2 = 96348
$ws.Range("E3").Value2 = 220787
$ws.Range("Q3").Value2 = 600749.0751519018
$ws.Range("R3").Value2 = 6613971.934424319
$ws.Range("D3").Value2 = 'VU'
$ws.Range("F3").Value2 = 'Knärot'
$ws.Range("G3").Value2 = 'Goodyera repens'
$ws.Range("H3").Value2 = '(L.) R. Br.'
$ws.Range("I3").Value2 = '10'
$ws.Range("J3").Value2 = 'plantor/tuvor'
$ws.Range("K3").Value2 = 'fullt utvecklade blad'

# Row 4  (was row 9)
$ws.Range("A4").Value2 = 110282856
$ws.Range("B4").Value2 = 89802
$ws.Range("E4").Value2 = 5420
$ws.Range("Q4").Value2 = 600677.6983460309
$ws.Range("R4").Value2 = 6613951.301940188
$ws.Range("D4").Value2 = 'LC'
$ws.Range("F4").Value2 = 'Grovticka'
$ws.Range("G4").Value2 = 'Phaeolus schweinitzii'
$ws.Range("H4").Value2 = '(Fr.) Pat.'
$ws.Range("I4").Value2 = '1'
$ws.Range("J4").Value2 = 'fruktkroppar'
$ws.Range("K4").Value2 = ''

# Row 5  (was row 7)
$ws.Range("A5").Value2 = 110282820
$ws.Range("B5").Value2 = 89425
$ws.Range("E5").Value2 = 5442
$ws.Range("Q5").Value2 = 600724.7123983201
$ws.Range("R5").Value2 = 6614086.574870056
$ws.Range("D5").Value2 = 'NT'
$ws.Range("F5").Value2 = 'Tallticka'
$ws.Range("G5").Value2 = 'Porodaedalea pini'
$ws.Range("H5").Value2 = '(Brot.) Murrill'
$ws.Range("I5").Value2 = '1'
$ws.Range("J5").Value2 = 'fruktkroppar'
$ws.Range("K5").Value2 = ''

# Row 6  (was row 2)
$ws.Range("A6").Value2 = 110282828
$ws.Range("B6").Value2 = 89425
$ws.Range("E6").Value2 = 5442
$ws.Range("Q6").Value2 = 600787.8656294679
$ws.Range("R6").Value2 = 6613904.709995793
$ws.Range("D6").Value2 = 'NT'
$ws.Range("F6").Value2 = 'Tallticka'
$ws.Range("G6").Value2 = 'Porodaedalea pini'
$ws.Range("H6").Value2 = '(Brot.) Murrill'
$ws.Range("I6").Value2 = '2'
$ws.Range("J6").Value2 = 'fruktkroppar'
$ws.Range("K6").Value2 = ''

# Row 7  (was row 6)
$ws.Range("A7").Value2 = 110282846
$ws.Range("B7").Value2 = 103288
$ws.Range("E7").Value2 = 221144
$ws.Range("Q7").Value2 = 600839.9318167433
$ws.Range("R7").Value2 = 6613983.990819811
$ws.Range("D7").Value2 = 'LC'
$ws.Range("F7").Value2 = 'Grönpyrola'
$ws.Range("G7").Value2 = 'Pyrola chlorantha'
$ws.Range("H7").Value2 = 'Sw.'
$ws.Range("I7").Value2 = '10'
$ws.Range("J7").Value2 = 'plantor/tuvor'
$ws.Range("K7").Value2 = 'blomning'

# Row 8  (was row 4)
$ws.Range("A8").Value2 = 110282836
$ws.Range("B8").Value2 = 89793
$ws.Range("E8").Value2 = 4217
$ws.Range("Q8").Value2 = 600805.3583702671
$ws.Range("R8").Value2 = 6613969.910894822
$ws.Range("D8").Value2 = 'LC'
$ws.Range("F8").Value2 = 'Blodticka'
$ws.Range("G8").Value2 = 'Meruliopsis taxicola'
$ws.Range("H8").Value2 = '(Pers.:Fr.) Bondartsev'
$ws.Range("I8").Value2 = '1'
$ws.Range("J8").Value2 = 'mycel'
$ws.Range("K8").Value2 = ''

# Row 9  (was row 3)
$ws.Range("A9").Value2 = 110282848
$ws.Range("B9").Value2 = 96348
$ws.Range("E9").Value2 = 220787
$ws.Range("Q9").Value2 = 600839.9318167433
$ws.Range("R9").Value2 = 6613983.990819811
$ws.Range("D9").Value2 = 'VU'
$ws.Range("F9").Value2 = 'Knärot'
$ws.Range("G9").Value2 = 'Goodyera repens'
$ws.Range("H9").Value2 = '(L.) R. Br.'
$ws.Range("I9").Value2 = '5'
$ws.Range("J9").Value2 = 'plantor/tuvor'
$ws.Range("K9").Value2 = 'fullt utvecklade blad'

# Cells that must become entirely empty (no explicit blank marker)
$ws.Range("L6").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AF8").ClearContents()
